$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "299.60"
Set-TextValue $ws.Range("E2") "1.86%"
Set-TextValue $ws.Range("D3") "32.22"
Set-TextValue $ws.Range("E3") "3.84%"
Set-TextValue $ws.Range("D4") "5.002"
Set-TextValue $ws.Range("E4") "1.51%"
Set-TextValue $ws.Range("D5") "0.07708"
Set-TextValue $ws.Range("E5") "4.95%"
Set-TextValue $ws.Range("D6") "2.253"
Set-TextValue $ws.Range("E6") "-2.01%"
Set-TextValue $ws.Range("D7") "7.929"
Set-TextValue $ws.Range("E7") "2.46%"
Set-TextValue $ws.Range("B8") "GateToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D8") "3.816"
Set-TextValue $ws.Range("E8") "1.79%"
Set-TextValue $ws.Range("B9") "MXToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9223"
Set-TextValue $ws.Range("E9") "1.51%"
Set-TextValue $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.09919"
Set-TextValue $ws.Range("E10") "24.62%"
Set-TextValue $ws.Range("B11") "WazirX"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1766"
Set-TextValue $ws.Range("E11") "4.63%"
Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.08398"
Set-TextValue $ws.Range("E12") "3.98%"
Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03301"
Set-TextValue $ws.Range("E13") "6.56%"
Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09835"
Set-TextValue $ws.Range("E14") "-2.47%"
Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001479"
Set-TextValue $ws.Range("E15") "-2.99%"
Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005654"
Set-TextValue $ws.Range("E16") "-3.57%"
Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.535"
Set-TextValue $ws.Range("E17") "1.45%"
Set-TextValue $ws.Range("D18") "2.196"
Set-TextValue $ws.Range("E18") "5.88%"
Set-TextValue $ws.Range("D19") "0.3373"
Set-TextValue $ws.Range("E19") "1.37%"
Set-TextValue $ws.Range("D20") "0.1336"
Set-TextValue $ws.Range("E20") "2.42%"
Set-TextValue $ws.Range("D21") "4.119"
Set-TextValue $ws.Range("E21") "3.67%"
Set-TextValue $ws.Range("E22") "-0.66%"
Set-TextValue $ws.Range("E23") "-0.44%"
Set-TextValue $ws.Range("D24") "0.001214"
Set-TextValue $ws.Range("E24") "0.32%"
Set-TextValue $ws.Range("D25") "0.004368"
Set-TextValue $ws.Range("E25") "-6.01%"
Set-TextValue $ws.Range("D26") "0.0001290"
Set-TextValue $ws.Range("E26") "-0.82%"
Set-TextValue $ws.Range("D27") "0.0003368"
Set-TextValue $ws.Range("E27") "-0.88%"
Set-TextValue $ws.Range("D39") "0.01708"
Set-TextValue $ws.Range("E39") "6.42%"
Set-TextValue $ws.Range("D40") "0.04656"
Set-TextValue $ws.Range("E40") "4.85%"
Set-TextValue $ws.Range("D41") "0.007638"
Set-TextValue $ws.Range("E41") "3.75%"
Set-TextValue $ws.Range("D42") "0.009755"
Set-TextValue $ws.Range("E42") "12.82%"
Set-TextValue $ws.Range("D43") "0.1394"
Set-TextValue $ws.Range("E43") "4.85%"
Set-TextValue $ws.Range("D44") "0.002085"
Set-TextValue $ws.Range("E44") "5.92%"
Set-TextValue $ws.Range("D45") "0.009709"
Set-TextValue $ws.Range("E45") "1.91%"
Set-TextValue $ws.Range("D46") "0.00006059"
Set-TextValue $ws.Range("E46") "1.85%"
Set-TextValue $ws.Range("D47") "0.00000000745"
Set-TextValue $ws.Range("E47") "-0.80%"
Set-TextValue $ws.Range("D48") "2.794"
Set-TextValue $ws.Range("E48") "24.68%"
Set-TextValue $ws.Range("D49") "0.001984"
Set-TextValue $ws.Range("E49") "-31.55%"
Set-TextValue $ws.Range("D50") "0.00002085"
Set-TextValue $ws.Range("E50") "-0.80%"
Set-TextValue $ws.Range("D51") "0.0001986"
Set-TextValue $ws.Range("E51") "-0.80%"
